$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the old "Decentral_*" rows (rows 12-20) - these are being replaced
# by the existing "nan_*" rows that follow them.
$ws.Range("A12:A20").EntireRow.Delete() | Out-Null

# After deletion, the rows that were 21-32 ("id_DK_nan_*") have shifted up
# to become rows 12-23, which matches the target layout. Now delete the
# trailing rows 24-32 that are no longer needed (dimension shrinks to B23).
$ws.Range("A24:A32").EntireRow.Delete() | Out-Null
